$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 0.1113992423925084
$ws.Range("B3").Value = 0.1254454475001058
$ws.Range("H3").Value = 0.2368446898926141
$ws.Range("B4").Value = 0.1490700599357073
$ws.Range("H4").Value = 0.2604693023282156
$ws.Range("B5").Value = 0.04494109595997438
$ws.Range("H5").Value = 0.1563403383524828
$ws.Range("B6").Value = 0.02783826860116842
$ws.Range("H6").Value = 0.1392375109936768
$ws.Range("B7").Value = 0.01649965745265265
$ws.Range("C7").Value = 0.001700272677293911
$ws.Range("D7").Value = 2.879625934947438
$ws.Range("E7").Value = 0.008800124277461713
$ws.Range("F7").Value = 0.01316146255509465
$ws.Range("G7").Value = 0.01983785235021041
$ws.Range("H7").Value = 0.127898899845161
$ws.Range("B8").Value = 0.01093051472373631
$ws.Range("H8").Value = 0.1223297571162447
$ws.Range("B9").Value = 0.0121749531417285
$ws.Range("H9").Value = 0.1235741955342369
$ws.Range("B10").Value = 0.01084051814073153
$ws.Range("C10").Value = 0.001058575450541831
$ws.Range("D10").Value = 2.026330108066947
$ws.Range("E10").Value = 0.003584510644042815
$ws.Range("F10").Value = 0.008760414580765916
$ws.Range("G10").Value = 0.01292062170069702
$ws.Range("H10").Value = 0.1222397605332399
$ws.Range("B11").Value = 0.02965607392827552
$ws.Range("H11").Value = 0.1410553163207839
$ws.Range("B12").Value = 0.05657811113143149
$ws.Range("H12").Value = 0.1679773535239399
$ws.Range("B13").Value = 0.07498615613373034
$ws.Range("H13").Value = 0.1863853985262387
$ws.Range("B14").Value = 0.08489125675829605
$ws.Range("H14").Value = 0.1962904991508044
$ws.Range("B15").Value = 0.09049349916045676
$ws.Range("H15").Value = 0.2018927415529652
$ws.Range("B16").Value = 0.09522265827235654
$ws.Range("H16").Value = 0.2066219006648649
$ws.Range("B17").Value = 0.09864656041426952
$ws.Range("H17").Value = 0.2100458028067779
$ws.Range("B18").Value = -0.1113992423925084
$ws.Range("C18").Value = 0.00788488177873389
$ws.Range("D18").Value = -23.33541603221571
$ws.Range("E18").Value = 0.02595550254954082
$ws.Range("F18").Value = -0.1269030552147783
$ws.Range("G18").Value = -0.09589542957023811
$ws.Range("B19").Value = 0.09987356293835122
$ws.Range("H19").Value = 0.2112728053308596
$ws.Range("B20").Value = 0.1021766204415719
$ws.Range("H20").Value = 0.2135758628340803
$ws.Range("B21").Value = 0.1059495378426912
$ws.Range("C21").Value = 0.006285832327742063
$ws.Range("D21").Value = 27.93579322505903
$ws.Range("E21").Value = 0.03486148792792863
$ws.Range("F21").Value = 0.09359442412155294
$ws.Range("G21").Value = 0.1183046515638294
$ws.Range("H21").Value = 0.2173487802351996
$ws.Range("B22").Value = 0.1101130185115464
$ws.Range("C22").Value = 0.006323449521945069
$ws.Range("D22").Value = 29.17957887082748
$ws.Range("E22").Value = 0.03421171292818506
$ws.Range("F22").Value = 0.09767520931049724
$ws.Range("G22").Value = 0.1225508277125959
$ws.Range("H22").Value = 0.2215122609040548
$ws.Range("B23").Value = 0.1158468035174292
$ws.Range("H23").Value = 0.2272460459099376
$ws.Range("B24").Value = 0.1176795602584984
$ws.Range("H24").Value = 0.2290788026510067
$ws.Range("B25").Value = 0.1204609444041502
$ws.Range("C25").Value = 0.00647845111161035
$ws.Range("D25").Value = 31.01524268844832
$ws.Range("E25").Value = 0.03271069084937731
$ws.Range("F25").Value = 0.1077210986172392
$ws.Range("G25").Value = 0.1332007901910614
$ws.Range("H25").Value = 0.2318601867966585
$ws.Range("B26").Value = 0.1229849557856659
$ws.Range("C26").Value = 0.006104691195107123
$ws.Range("D26").Value = 31.84669465832229
$ws.Range("E26").Value = 0.03420637287696614
$ws.Range("F26").Value = 0.1109793702228082
$ws.Range("G26").Value = 0.134990541348524
$ws.Range("H26").Value = 0.2343841981781743
$ws.Range("B27").Value = 0.1293952131760527
$ws.Range("C27").Value = 0.005807337085421548
$ws.Range("D27").Value = 32.10128645339148
$ws.Range("E27").Value = 0.03514030670067056
$ws.Range("F27").Value = 0.117993919839423
$ws.Range("G27").Value = 0.1407965065126831
$ws.Range("H27").Value = 0.2407944555685611
$ws.Range("B28").Value = 0.1311854579586333
$ws.Range("C28").Value = 0.005934526194959007
$ws.Range("D28").Value = 31.39143688929453
$ws.Range("E28").Value = 0.06214760782429639
$ws.Range("F28").Value = 0.1195352258102435
$ws.Range("G28").Value = 0.142835690107023
$ws.Range("H28").Value = 0.2425847003511417
$ws.Range("B29").Value = 0.01316488972226798
$ws.Range("C29").Value = 0.0007613288916396197
$ws.Range("D29").Value = 2.496460847919775
$ws.Range("E29").Value = 0.003316006460320914
$ws.Range("F29").Value = 0.01166732989672211
$ws.Range("G29").Value = 0.01466244954781373
$ws.Range("H29").Value = 0.1245641321147764
